$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1054.4286  # H15
$ws.Cells.Item(15, 9).Value = 1054.4286  # I15
$ws.Cells.Item(15, 11).Value = 3163.2858  # K15
$ws.Cells.Item(15, 13).Value = -2994.2858  # M15
$ws.Cells.Item(17, 8).Value = 766142.75  # H17
$ws.Cells.Item(17, 10).Value = 766142.75  # J17
$ws.Cells.Item(17, 12).Value = 2298428.25  # L17
$ws.Cells.Item(17, 14).Value = -2298764.25  # N17
$ws.Cells.Item(132, 8).Value = 2345.4253  # H132
$ws.Cells.Item(132, 9).Value = 1946.0509  # I132
$ws.Cells.Item(132, 10).Value = 3186.9644  # J132
$ws.Cells.Item(132, 11).Value = 5838.1527  # K132
$ws.Cells.Item(132, 12).Value = 9560.893199999999  # L132
$ws.Cells.Item(132, 13).Value = -3308.1527  # M132
$ws.Cells.Item(132, 14).Value = -14620.8932  # N132
$ws.Cells.Item(137, 8).Value = 3302.1147  # H137
$ws.Cells.Item(137, 9).Value = 1422.2858  # I137
$ws.Cells.Item(137, 10).Value = 4897.121  # J137
$ws.Cells.Item(137, 11).Value = 4266.857400000001  # K137
$ws.Cells.Item(137, 12).Value = 14691.363  # L137
$ws.Cells.Item(137, 13).Value = -1716.857400000001  # M137
$ws.Cells.Item(137, 14).Value = -19791.363  # N137
$ws.Cells.Item(138, 8).Value = 2718.7446  # H138
$ws.Cells.Item(138, 9).Value = 1706.5294  # I138
$ws.Cells.Item(138, 10).Value = 3292.3333  # J138
$ws.Cells.Item(138, 11).Value = 5119.5882  # K138
$ws.Cells.Item(138, 12).Value = 9876.999899999999  # L138
$ws.Cells.Item(138, 13).Value = 20.41179999999986  # M138
$ws.Cells.Item(138, 14).Value = -20156.9999  # N138
$ws.Cells.Item(141, 8).Value = 2697.2952  # H141
$ws.Cells.Item(141, 9).Value = 932.6786  # I141
$ws.Cells.Item(141, 11).Value = 2798.0358  # K141
$ws.Cells.Item(141, 13).Value = 2381.9642  # M141

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4837.152  # H32
$ws.Cells.Item(32, 9).Value = 3834.2239  # I32
$ws.Cells.Item(32, 10).Value = 10436.833  # J32
$ws.Cells.Item(32, 11).Value = 3834.2239  # K32
$ws.Cells.Item(32, 12).Value = 10436.833  # L32
$ws.Cells.Item(32, 13).Value = -3547.2239  # M32
$ws.Cells.Item(32, 14).Value = -11010.833  # N32
$ws.Cells.Item(61, 8).Value = 1580.9231  # H61
$ws.Cells.Item(61, 9).Value = 1401.8334  # I61
$ws.Cells.Item(61, 10).Value = 2177.889  # J61
$ws.Cells.Item(61, 11).Value = 1401.8334  # K61
$ws.Cells.Item(61, 12).Value = 2177.889  # L61
$ws.Cells.Item(61, 13).Value = -1189.8334  # M61
$ws.Cells.Item(61, 14).Value = -2601.889  # N61
$ws.Cells.Item(74, 8).Value = 4923.68  # H74
$ws.Cells.Item(74, 9).Value = 876.2727  # I74
$ws.Cells.Item(74, 10).Value = 34604.668  # J74
$ws.Cells.Item(74, 11).Value = 876.2727  # K74
$ws.Cells.Item(74, 12).Value = 34604.668  # L74
$ws.Cells.Item(74, 13).Value = -2.272699999999986  # M74
$ws.Cells.Item(74, 14).Value = -36352.668  # N74
$ws.Cells.Item(75, 8).Value = 31586.5  # H75
$ws.Cells.Item(75, 10).Value = 31586.5  # J75
$ws.Cells.Item(75, 12).Value = 31586.5  # L75
$ws.Cells.Item(75, 14).Value = -33334.5  # N75
$ws.Cells.Item(77, 8).Value = 4923.68  # H77
$ws.Cells.Item(77, 9).Value = 876.2727  # I77
$ws.Cells.Item(77, 10).Value = 34604.668  # J77
$ws.Cells.Item(77, 11).Value = 4381.363499999999  # K77
$ws.Cells.Item(77, 12).Value = 173023.34  # L77
$ws.Cells.Item(77, 13).Value = -13.36349999999948  # M77
$ws.Cells.Item(77, 14).Value = -181759.34  # N77
$ws.Cells.Item(78, 8).Value = 31586.5  # H78
$ws.Cells.Item(78, 10).Value = 31586.5  # J78
$ws.Cells.Item(78, 12).Value = 94759.5  # L78
$ws.Cells.Item(78, 14).Value = -103495.5  # N78
$ws.Cells.Item(88, 8).Value = 3858.9  # H88
$ws.Cells.Item(88, 9).Value = 3254.5715  # I88
$ws.Cells.Item(88, 10).Value = 5269  # J88
$ws.Cells.Item(88, 11).Value = 3254.5715  # K88
$ws.Cells.Item(88, 12).Value = 5269  # L88
$ws.Cells.Item(88, 13).Value = -2848.5715  # M88
$ws.Cells.Item(88, 14).Value = -6081  # N88
$ws.Cells.Item(91, 8).Value = 3858.9  # H91
$ws.Cells.Item(91, 9).Value = 3254.5715  # I91
$ws.Cells.Item(91, 10).Value = 5269  # J91
$ws.Cells.Item(91, 11).Value = 3254.5715  # K91
$ws.Cells.Item(91, 12).Value = 5269  # L91
$ws.Cells.Item(91, 13).Value = -1850.5715  # M91
$ws.Cells.Item(91, 14).Value = -8077  # N91
$ws.Cells.Item(136, 8).Value = 1580.9231  # H136
$ws.Cells.Item(136, 9).Value = 1401.8334  # I136
$ws.Cells.Item(136, 10).Value = 2177.889  # J136
$ws.Cells.Item(136, 11).Value = 4205.5002  # K136
$ws.Cells.Item(136, 12).Value = 6533.667  # L136
$ws.Cells.Item(136, 13).Value = -1655.5002  # M136
$ws.Cells.Item(136, 14).Value = -11633.667  # N136

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3172.4666  # H86
$ws.Cells.Item(86, 9).Value = 3961.2  # I86
$ws.Cells.Item(86, 11).Value = 3961.2  # K86
$ws.Cells.Item(86, 13).Value = -2838.2  # M86
$ws.Cells.Item(89, 8).Value = 3172.4666  # H89
$ws.Cells.Item(89, 9).Value = 3961.2  # I89
$ws.Cells.Item(89, 11).Value = 19806  # K89
$ws.Cells.Item(89, 13).Value = -14190  # M89
$ws.Cells.Item(134, 8).Value = 1127.8235  # H134
$ws.Cells.Item(134, 9).Value = 1027.4814  # I134
$ws.Cells.Item(134, 10).Value = 1514.8572  # J134
$ws.Cells.Item(134, 11).Value = 3082.4442  # K134
$ws.Cells.Item(134, 12).Value = 4544.571599999999  # L134
$ws.Cells.Item(134, 13).Value = -547.4441999999999  # M134
$ws.Cells.Item(134, 14).Value = -9614.571599999999  # N134

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 27594.873  # H31
$ws.Cells.Item(31, 9).Value = 2951.389  # I31
$ws.Cells.Item(31, 10).Value = 69180.75  # J31
$ws.Cells.Item(31, 11).Value = 2951.389  # K31
$ws.Cells.Item(31, 12).Value = 69180.75  # L31
$ws.Cells.Item(31, 13).Value = -2656.389  # M31
$ws.Cells.Item(31, 14).Value = -69770.75  # N31
$ws.Cells.Item(34, 8).Value = 27594.873  # H34
$ws.Cells.Item(34, 9).Value = 2951.389  # I34
$ws.Cells.Item(34, 10).Value = 69180.75  # J34
$ws.Cells.Item(34, 11).Value = 2951.389  # K34
$ws.Cells.Item(34, 12).Value = 69180.75  # L34
$ws.Cells.Item(34, 13).Value = -2749.389  # M34
$ws.Cells.Item(34, 14).Value = -69584.75  # N34
$ws.Cells.Item(58, 8).Value = 854.7  # H58
$ws.Cells.Item(58, 9).Value = 747  # I58
$ws.Cells.Item(58, 10).Value = 1393.2  # J58
$ws.Cells.Item(58, 11).Value = 747  # K58
$ws.Cells.Item(58, 12).Value = 1393.2  # L58
$ws.Cells.Item(58, 13).Value = -544  # M58
$ws.Cells.Item(58, 14).Value = -1799.2  # N58
$ws.Cells.Item(62, 8).Value = 4141.6387  # H62
$ws.Cells.Item(62, 9).Value = 4950.227  # I62
$ws.Cells.Item(62, 10).Value = 2871  # J62
$ws.Cells.Item(62, 11).Value = 4950.227  # K62
$ws.Cells.Item(62, 12).Value = 2871  # L62
$ws.Cells.Item(62, 13).Value = -4326.227  # M62
$ws.Cells.Item(62, 14).Value = -4119  # N62
$ws.Cells.Item(65, 8).Value = 4141.6387  # H65
$ws.Cells.Item(65, 9).Value = 4950.227  # I65
$ws.Cells.Item(65, 10).Value = 2871  # J65
$ws.Cells.Item(65, 11).Value = 24751.135  # K65
$ws.Cells.Item(65, 12).Value = 14355  # L65
$ws.Cells.Item(65, 13).Value = -21631.135  # M65
$ws.Cells.Item(65, 14).Value = -20595  # N65
$ws.Cells.Item(132, 8).Value = 13336432  # H132
$ws.Cells.Item(132, 9).Value = 19234316  # I132
$ws.Cells.Item(132, 10).Value = 2083.8262  # J132
$ws.Cells.Item(132, 11).Value = 57702948  # K132
$ws.Cells.Item(132, 12).Value = 6251.4786  # L132
$ws.Cells.Item(132, 13).Value = -57700418  # M132
$ws.Cells.Item(132, 14).Value = -11311.4786  # N132
$ws.Cells.Item(136, 8).Value = 854.7  # H136
$ws.Cells.Item(136, 9).Value = 747  # I136
$ws.Cells.Item(136, 10).Value = 1393.2  # J136
$ws.Cells.Item(136, 11).Value = 2241  # K136
$ws.Cells.Item(136, 12).Value = 4179.6  # L136
$ws.Cells.Item(136, 13).Value = 309  # M136
$ws.Cells.Item(136, 14).Value = -9279.6  # N136

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(123, 8).Value = 26595.6  # H123
$ws.Cells.Item(123, 10).Value = 26595.6  # J123
$ws.Cells.Item(123, 12).Value = 26595.6  # L123
$ws.Cells.Item(123, 14).Value = -31495.6  # N123
$ws.Cells.Item(132, 8).Value = 5572.613  # H132
$ws.Cells.Item(132, 9).Value = 7537.579  # I132
$ws.Cells.Item(132, 10).Value = 2461.4167  # J132
$ws.Cells.Item(132, 11).Value = 22612.737  # K132
$ws.Cells.Item(132, 12).Value = 7384.250100000001  # L132
$ws.Cells.Item(132, 13).Value = -20082.737  # M132
$ws.Cells.Item(132, 14).Value = -12444.2501  # N132

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2231.2222  # H46
$ws.Cells.Item(46, 9).Value = 2322.625  # I46
$ws.Cells.Item(46, 10).Value = 1500  # J46
$ws.Cells.Item(46, 11).Value = 2322.625  # K46
$ws.Cells.Item(46, 12).Value = 1500  # L46
$ws.Cells.Item(46, 13).Value = -2134.625  # M46
$ws.Cells.Item(46, 14).Value = -1876  # N46
$ws.Cells.Item(136, 8).Value = 2799.6897  # H136
$ws.Cells.Item(136, 9).Value = 981.3913  # I136
$ws.Cells.Item(136, 10).Value = 9769.833000000001  # J136
$ws.Cells.Item(136, 11).Value = 2944.1739  # K136
$ws.Cells.Item(136, 12).Value = 29309.499  # L136
$ws.Cells.Item(136, 13).Value = -394.1738999999998  # M136
$ws.Cells.Item(136, 14).Value = -34409.499  # N136
$ws.Cells.Item(140, 8).Value = 40500  # H140
$ws.Cells.Item(140, 10).Value = 40500  # J140
$ws.Cells.Item(140, 12).Value = 40500  # L140
$ws.Cells.Item(140, 14).Value = -50860  # N140

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 600  # H107
$ws.Cells.Item(107, 9).Value = 600  # I107
$ws.Cells.Item(107, 10).Value = 0  # J107
$ws.Cells.Item(107, 11).Value = 1800  # K107
$ws.Cells.Item(107, 12).Value = 0  # L107
$ws.Cells.Item(107, 13).Value = 120  # M107
$ws.Cells.Item(107, 14).ClearContents()  # N107
$ws.Cells.Item(132, 8).Value = 2789.4517  # H132
$ws.Cells.Item(132, 9).Value = 3505.244  # I132
$ws.Cells.Item(132, 10).Value = 1391.9524  # J132
$ws.Cells.Item(132, 11).Value = 10515.732  # K132
$ws.Cells.Item(132, 12).Value = 4175.857199999999  # L132
$ws.Cells.Item(132, 13).Value = -7985.732  # M132
$ws.Cells.Item(132, 14).Value = -9235.857199999999  # N132
$ws.Cells.Item(136, 8).Value = 2463.761  # H136
$ws.Cells.Item(136, 9).Value = 3369.1875  # I136
$ws.Cells.Item(136, 10).Value = 1476.0227  # J136
$ws.Cells.Item(136, 11).Value = 10107.5625  # K136
$ws.Cells.Item(136, 12).Value = 4428.0681  # L136
$ws.Cells.Item(136, 13).Value = -7557.5625  # M136
$ws.Cells.Item(136, 14).Value = -9528.0681  # N136
